$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: box-plots / histogram paragraph - rewrite text
# -----------------------------------------------------------------
$old1 = "In the univariate case, box-plots do provide some information that the histogram does not (at least, not explicitly). That is, it typically provides the median, 25th and 75th percentile, min/max that is not an outlier and explicitly separates the points that are considered outliers. This can all be ""eyeballed"" from the histogram and may be better to be eyeballed in the case of outliers."
$new1 = "In the univariate case, box-plots do provide some information that the histogram does not. That is, it typically provides the median, 25th and 75th percentile, min/max that is not an outlier and explicitly separates the points that are considered outliers. In box plots we can better visualize the outliers with respect to the inferential statistics of the feature. Histograms only give a meaure of desity of the feature values"
$found = $d.Content.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# -----------------------------------------------------------------
# Change 2: "would the results ... would be applied" rewording
# (search includes the preceding word "how" so both now-orphaned
#  proofErr gramStart/gramEnd markers get absorbed/removed cleanly)
# -----------------------------------------------------------------
$old2 = "how would the results would be applied in real word."
$new2 = "how the results would  would be applied in real word."
$found = $d.Content.Find.Execute($old2, $false, $false, $false, $false, $false, $true, 1, $false, $new2, 2)

# -----------------------------------------------------------------
# Change 3: merge the split runs in the hypothesis-testing paragraphs
# (self-replace forces the runtime to coalesce runs with identical
#  formatting into a single run, matching the target XML)
# -----------------------------------------------------------------
$same3a = " choose a suitable statistical test and statistics used to reject the null hypothesis and choose a critical region for the statistics to lie in that is extreme enough for the null hypothesis to be rejected (p-value)"
$found = $d.Content.Find.Execute($same3a, $false, $false, $false, $false, $false, $true, 1, $false, $same3a, 2)

$same3b = "We then calculate the observed test statistics from the data and check whether it lies in the critical region. There are multiple test we performed based on the nature of the problem and features of our dataset."
$found = $d.Content.Find.Execute($same3b, $false, $false, $false, $false, $false, $true, 1, $false, $same3b, 2)

# -----------------------------------------------------------------
# Change 4: apply Arial font formatting to the "Types of distribution" paragraph
# -----------------------------------------------------------------
$found4 = $d.Content.Find.Execute("Types of distribution that are non-Gaussian or non-log normal are the skewed distributions, discrete distributions and binomial distribution.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found4) {
  $rng4 = $d.Content
  $rng4.Start = $rng4.Start
  $p4 = $d.Paragraphs(25)
  $p4.Range.Font.Name = "Arial"
}

# -----------------------------------------------------------------
# Change 6: extend "central tendency than mean." paragraph with new sentence
# -----------------------------------------------------------------
$rngA = $d.Content
$findA = $rngA.Find
$findA.ClearFormatting()
$findA.Text = "central tendency than mean. "
$foundA = $findA.Execute()
if ($foundA) {
  $rngA.Collapse(0)
  $rngA.InsertAfter(" . Another time when we usually prefer the median over the mean (or mode) is when our data is skewed")
}

# -----------------------------------------------------------------
# Change 7: drop the leading ". " before "A likelihood function"
# -----------------------------------------------------------------
$old7 = ". A likelihood function"
$new7 = "A likelihood function"
$found = $d.Content.Find.Execute($old7, $false, $false, $false, $false, $false, $true, 1, $false, $new7, 2)

# -----------------------------------------------------------------
# Change 8: extend the likelihood-function paragraph with further text,
# then relocate the _GoBack bookmark into the middle of the new text
# (splitting "T" | "he" of "The median"), leaving the two trailing
# ListParagraph paragraphs empty.
# -----------------------------------------------------------------
$rngB = $d.Content
$findB = $rngB.Find
$findB.ClearFormatting()
$findB.Text = "how well the data summarizes these parameters."
$foundB = $findB.Execute()
if ($foundB) {
  $rngB.Collapse(0)
  $rngB.InsertAfter(" Maximum likelihood estimation is a method that determines values for the parameters of a model. The parameter values are found such that they maximise the likelihood that the process described by the model produced the data that were actually observed. The median best retains this position and is not as strongly influenced by the skewed values.")
}

$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$rngC = $d.Content
$findC = $rngC.Find
$findC.ClearFormatting()
$findC.Text = "The median best retains"
$foundC = $findC.Execute()
if ($foundC) {
  $splitPos = $rngC.Start + 1
  $bmRange = $d.Range($splitPos, $splitPos)
  $d.Bookmarks.Add("_GoBack", $bmRange)
}

Write-Output "done"
